$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = 'stomach compression for women'
$ws.Range("A2").Value = 'stomach exercises'
$ws.Range("A3").Value = 'stomach workout'
$ws.Range("A4").Value = 'strained muscle'
$ws.Range("A5").Value = 'stretch calves'
$ws.Range("A6").Value = 'stretch capri pants'
$ws.Range("A7").Value = 'stretch hamstring'
$ws.Range("A8").Value = 'stretch hiking pants heavy'
$ws.Range("A9").Value = 'stretch legging'
$ws.Range("A10").Value = 'stretch legs'
$ws.Range("A11").Value = 'stretch muscle'
$ws.Range("A12").Value = 'stretch muscles'
$ws.Range("A13").Value = 'stretch quad'
$ws.Range("A14").Value = 'stretch rings for women'
$ws.Range("A15").Value = 'stretch tight'
$ws.Range("A16").Value = 'stretch tights'
$ws.Range("A17").Value = 'stretch yoga leggings'
$ws.Range("A18").Value = 'stretching gear'
$ws.Range("A19").Value = 'stretching pants for women'
$ws.Range("A20").Value = 'stretching pants women'
$ws.Range("A21").Value = 'stretchy ballet flats'
$ws.Range("A22").Value = 'stretchy tights for women'
$ws.Range("A23").Value = 'string art bike'
$ws.Range("A24").Value = 'stripe capris'
$ws.Range("A25").Value = 'stripes thigh highs'
$ws.Range("A26").Value = 'stroke flash'
$ws.Range("A27").Value = 'stroke recovery equipment'
$ws.Range("A28").Value = 'strong joints'
$ws.Range("A29").Value = 'strong sleeves'
$ws.Range("A30").Value = 'style and company pants for women'
$ws.Range("A31").Value = 'subsports compression'
$ws.Range("A32").Value = 'summer breeches'
$ws.Range("A33").Value = 'summer capri leggings'
$ws.Range("A34").Value = 'summer capri leggings for women'
$ws.Range("A35").Value = 'summer capris for women'
$ws.Range("A36").Value = 'summer clothes for women over 50'
$ws.Range("A37").Value = 'summer leggings for women'
$ws.Range("A38").Value = 'summer leggings for women capri'
$ws.Range("A39").Value = 'summer riding pants'
$ws.Range("A40").Value = 'summer waist slimmer'
$ws.Range("A41").Value = 'summer workout leggings'
$ws.Range("A42").Value = 'sun leggings'
$ws.Range("A43").Value = 'sun squad slip and slide'
$ws.Range("A44").Value = 'sunday bikes'
$ws.Range("A45").Value = 'super comfy leggings'
$ws.Range("A46").Value = 'super compression leggings'
$ws.Range("A47").Value = 'super high rise leggings'
$ws.Range("A48").Value = 'super high waisted black leggings'
$ws.Range("A49").Value = 'super joint support'
$ws.Range("A50").Value = 'super man leggings'
$ws.Range("A51").Value = 'super opaque tights'
$ws.Range("A52").Value = 'super rugby'
$ws.Range("A53").Value = 'super soft black leggings'
$ws.Range("A54").Value = 'super tight leggings'
$ws.Range("A55").Value = 'suport back brace'
$ws.Range("A56").Value = 'suport hose'
$ws.Range("A57").Value = 'support black tight'
$ws.Range("A58").Value = 'support capri'
$ws.Range("A59").Value = 'support for knees'
$ws.Range("A60").Value = 'support for spine'
$ws.Range("A61").Value = 'support gear'
$ws.Range("A62").Value = 'support hoops for garden fabric'
$ws.Range("A63").Value = 'support joint'
$ws.Range("A64").Value = 'support knee'
$ws.Range("A65").Value = 'support knee for women'
$ws.Range("A66").Value = 'support knee highs'
$ws.Range("A67").Value = 'support knee highs black'
$ws.Range("A68").Value = 'support knee highs for women'
$ws.Range("A69").Value = 'support knee highs for women plus size'
$ws.Range("A70").Value = 'support knee highs plus size'
$ws.Range("A71").Value = 'support knee highs women'
$ws.Range("A72").Value = 'support leggings'
$ws.Range("A73").Value = 'support leggings for women'
$ws.Range("A74").Value = 'support pants for men'
$ws.Range("A75").Value = 'support panty hose'
$ws.Range("A76").Value = 'support pantyhose for women black'
$ws.Range("A77").Value = 'support pany hose for women'
$ws.Range("A78").Value = 'support plus stockings'
$ws.Range("A79").Value = 'support shorts'
$ws.Range("A80").Value = 'support stick for walking'
$ws.Range("A81").Value = 'support tape'
$ws.Range("A82").Value = 'support tight'
$ws.Range("A83").Value = 'support tights'
$ws.Range("A84").Value = 'support tights for women'
$ws.Range("A85").Value = 'support tights plus size'
$ws.Range("A86").Value = 'support wear for women'
$ws.Range("A87").Value = 'support women'
$ws.Range("A88").Value = 'support yoga pants for women'
$ws.Range("A89").Value = 'supportive knee brace for women'
$ws.Range("A90").Value = 'surgery recovery pants'
$ws.Range("A91").Value = 'surgical clothing'
$ws.Range("A92").Value = 'surgical pants'
$ws.Range("A93").Value = 'surgical wear'
$ws.Range("A94").Value = 'survival bike gear'
$ws.Range("A95").Value = 'survival clothing women'
$ws.Range("A96").Value = 'survival gear vest'
$ws.Range("A97").Value = 'survival generator'
$ws.Range("A98").Value = 'survival power pot'
$ws.Range("A99").Value = 'sweat compression for women'
$ws.Range("A100").Value = 'sweat neoprene pants'
